$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 591.75
$ws.Range("I6").Value = 412.625
$ws.Range("J6").Value = 950
$ws.Range("K6").Value = 1237.875
$ws.Range("L6").Value = 2850
$ws.Range("M6").Value = -1125.875
$ws.Range("N6").Value = -3074
$ws.Range("H8").Value = 6038.4
$ws.Range("I8").Value = 48
$ws.Range("J8").Value = 30000
$ws.Range("K8").Value = 144
$ws.Range("L8").Value = 90000
$ws.Range("M8").Value = -5
$ws.Range("N8").Value = -90278
$ws.Range("H38").Value = 905.125
$ws.Range("I38").Value = 86.916664
$ws.Range("J38").Value = 3359.75
$ws.Range("K38").Value = 260.749992
$ws.Range("L38").Value = 10079.25
$ws.Range("M38").Value = 111.250008
$ws.Range("N38").Value = -10823.25
$ws.Range("H132").Value = 1560.6177
$ws.Range("I132").Value = 1430.7742
$ws.Range("J132").Value = 2902.3333
$ws.Range("K132").Value = 4292.3226
$ws.Range("L132").Value = 8706.999899999999
$ws.Range("M132").Value = -1762.3226
$ws.Range("N132").Value = -13766.9999
$ws.Range("H137").Value = 2066.5278
$ws.Range("I137").Value = 1960.4546
$ws.Range("J137").Value = 2233.2144
$ws.Range("K137").Value = 5881.3638
$ws.Range("L137").Value = 6699.6432
$ws.Range("M137").Value = -3331.3638
$ws.Range("N137").Value = -11799.6432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15824.027
$ws.Range("I32").Value = 12312.478
$ws.Range("J32").Value = 62878.8
$ws.Range("K32").Value = 12312.478
$ws.Range("L32").Value = 62878.8
$ws.Range("M32").Value = -12025.478
$ws.Range("N32").Value = -63452.8
$ws.Range("H88").Value = 3750
$ws.Range("I88").Value = 2500
$ws.Range("K88").Value = 2500
$ws.Range("M88").Value = -2094
$ws.Range("H91").Value = 3750
$ws.Range("I91").Value = 2500
$ws.Range("K91").Value = 2500
$ws.Range("M91").Value = -1096
$ws.Range("H124").Value = 42000
$ws.Range("J124").Value = 42000
$ws.Range("L124").Value = 42000
$ws.Range("N124").Value = -51820
$ws.Range("H125").Value = 137515
$ws.Range("J125").Value = 137515
$ws.Range("L125").Value = 137515
$ws.Range("N125").Value = -147355
$ws.Range("H139").Value = 84985.664
$ws.Range("J139").Value = 84985.664
$ws.Range("L139").Value = 84985.664
$ws.Range("N139").Value = -95265.664

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 258823.67
$ws.Range("I26").Value = 258823.67
$ws.Range("K26").Value = 258823.67
$ws.Range("M26").Value = -258531.67
$ws.Range("H86").Value = 85316
$ws.Range("I86").Value = 2163
$ws.Range("J86").Value = 144711
$ws.Range("K86").Value = 2163
$ws.Range("L86").Value = 144711
$ws.Range("M86").Value = -1040
$ws.Range("N86").Value = -146957
$ws.Range("H89").Value = 85316
$ws.Range("I89").Value = 2163
$ws.Range("J89").Value = 144711
$ws.Range("K89").Value = 10815
$ws.Range("L89").Value = 723555
$ws.Range("M89").Value = -5199
$ws.Range("N89").Value = -734787
$ws.Range("H134").Value = 502988.06
$ws.Range("I134").Value = 771580.9399999999
$ws.Range("J134").Value = 4172.7144
$ws.Range("K134").Value = 2314742.82
$ws.Range("L134").Value = 12518.1432
$ws.Range("M134").Value = -2312207.82
$ws.Range("N134").Value = -17588.1432

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9332.833000000001
$ws.Range("J4").Value = 9332.833000000001
$ws.Range("L4").Value = 9332.833000000001
$ws.Range("N4").Value = -9556.833000000001
$ws.Range("H31").Value = 3605.2222
$ws.Range("I31").Value = 2547.0557
$ws.Range("K31").Value = 2547.0557
$ws.Range("M31").Value = -2252.0557
$ws.Range("H34").Value = 3605.2222
$ws.Range("I34").Value = 2547.0557
$ws.Range("K34").Value = 2547.0557
$ws.Range("M34").Value = -2345.0557
$ws.Range("H94").Value = 2054.2307
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 2100.4167
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 2100.4167
$ws.Range("M94").Value = -1049
$ws.Range("N94").Value = -3002.4167
$ws.Range("H99").Value = 3177.8096
$ws.Range("I99").Value = 3070.875
$ws.Range("J99").Value = 3520
$ws.Range("K99").Value = 3070.875
$ws.Range("L99").Value = 3520
$ws.Range("M99").Value = -1572.875
$ws.Range("N99").Value = -6516
$ws.Range("H126").Value = 3177.8096
$ws.Range("I126").Value = 3070.875
$ws.Range("J126").Value = 3520
$ws.Range("K126").Value = 9212.625
$ws.Range("L126").Value = 10560
$ws.Range("M126").Value = -6742.625
$ws.Range("N126").Value = -15500
$ws.Range("H132").Value = 2230.9333
$ws.Range("I132").Value = 1703.9131
$ws.Range("J132").Value = 3962.5715
$ws.Range("K132").Value = 5111.7393
$ws.Range("L132").Value = 11887.7145
$ws.Range("M132").Value = -2581.7393
$ws.Range("N132").Value = -16947.7145
$ws.Range("H134").Value = 1774.3182
$ws.Range("I134").Value = 1739.6875
$ws.Range("J134").Value = 1866.6666
$ws.Range("K134").Value = 5219.0625
$ws.Range("L134").Value = 5599.9998
$ws.Range("M134").Value = -2684.0625
$ws.Range("N134").Value = -10669.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 107.14286
$ws.Range("I7").Value = 108.333336
$ws.Range("K7").Value = 325.000008
$ws.Range("M7").Value = -213.000008
$ws.Range("H12").Value = 1136751.1
$ws.Range("I12").Value = 108
$ws.Range("J12").Value = 1486487.5
$ws.Range("K12").Value = 324
$ws.Range("L12").Value = 4459462.5
$ws.Range("M12").Value = -151
$ws.Range("N12").Value = -4459808.5
$ws.Range("H41").Value = 804.1667
$ws.Range("J41").Value = 804.1667
$ws.Range("L41").Value = 2412.5001
$ws.Range("N41").Value = -3088.5001
$ws.Range("H46").Value = 3469.1875
$ws.Range("I46").Value = 233.33333
$ws.Range("J46").Value = 4215.923
$ws.Range("K46").Value = 699.99999
$ws.Range("L46").Value = 12647.769
$ws.Range("M46").Value = -608.99999
$ws.Range("N46").Value = -12829.769
$ws.Range("H116").Value = 1734.3846
$ws.Range("I116").Value = 738.3333
$ws.Range("J116").Value = 2033.2
$ws.Range("K116").Value = 2214.9999
$ws.Range("L116").Value = 6099.6
$ws.Range("M116").Value = 1227.0001
$ws.Range("N116").Value = -12983.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21064
$ws.Range("H95").Value = 10000000
$ws.Range("J95").Value = 10000000
$ws.Range("L95").Value = 10000000
$ws.Range("N95").Value = -10005492
$ws.Range("H132").Value = 2587.5386
$ws.Range("I132").Value = 1737.7778
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 5213.3334
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -2683.3334
$ws.Range("N132").Value = -18558.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9956.462
$ws.Range("J2").Value = 9956.462
$ws.Range("L2").Value = 9956.462
$ws.Range("N2").Value = -10180.462
$ws.Range("H7").Value = 4083.8333
$ws.Range("I7").Value = 5168
$ws.Range("J7").Value = 2999.6667
$ws.Range("K7").Value = 5168
$ws.Range("L7").Value = 2999.6667
$ws.Range("M7").Value = -5056
$ws.Range("N7").Value = -3223.6667
$ws.Range("H122").Value = 5388.2256
$ws.Range("I122").Value = 5494.5654
$ws.Range("J122").Value = 5082.5
$ws.Range("K122").Value = 16483.6962
$ws.Range("L122").Value = 15247.5
$ws.Range("M122").Value = -14033.6962
$ws.Range("N122").Value = -20147.5
$ws.Range("H126").Value = 4083.8333
$ws.Range("I126").Value = 5168
$ws.Range("J126").Value = 2999.6667
$ws.Range("K126").Value = 15504
$ws.Range("L126").Value = 8999.000100000001
$ws.Range("M126").Value = -13034
$ws.Range("N126").Value = -13939.0001
$ws.Range("H132").Value = 5139.3945
$ws.Range("I132").Value = 5686.077
$ws.Range("J132").Value = 3954.9167
$ws.Range("K132").Value = 17058.231
$ws.Range("L132").Value = 11864.7501
$ws.Range("M132").Value = -14528.231
$ws.Range("N132").Value = -16924.7501
$ws.Range("H140").Value = 51909.332
$ws.Range("J140").Value = 51909.332
$ws.Range("L140").Value = 51909.332
$ws.Range("N140").Value = -62269.332

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 15270.214
$ws.Range("I100").Value = 33646.668
$ws.Range("K100").Value = 67293.336
$ws.Range("M100").Value = -66752.336
$ws.Range("H109").Value = 40427
$ws.Range("J109").Value = 40427
$ws.Range("L109").Value = 40427
$ws.Range("N109").Value = -43201
$ws.Range("H122").Value = 6143.4287
$ws.Range("I122").Value = 4668
$ws.Range("J122").Value = 7250
$ws.Range("K122").Value = 14004
$ws.Range("L122").Value = 21750
$ws.Range("M122").Value = -11554
$ws.Range("N122").Value = -26650
